$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.81105205992511
$ws.Range("L2").Value = 0.804975443687201

$ws.Range("B3").Value = 0.79693552228442
$ws.Range("L3").Value = 0.8034001106575

$ws.Range("B4").Value = 0.757433918108113
$ws.Range("I4").Value = 0.742113769980062
$ws.Range("L4").Value = 0.847170664064891

$ws.Range("B5").Value = 0.711494502021427
$ws.Range("L5").Value = 0.729127889064831

$ws.Range("B6").Value = 0.711234316950339
$ws.Range("L6").Value = 0.810580370309662

$ws.Range("B7").Value = 0.710919543617656
$ws.Range("L7").Value = 0.528025667438765

$ws.Range("B8").Value = 0.698333620716356
$ws.Range("D8").Value = 0.609435949401637
$ws.Range("J8").Value = 0.615967974842448
$ws.Range("K8").Value = 0.682093186257981
$ws.Range("L8").Value = 0.726511394674859

$ws.Range("B9").Value = 0.695603247634989
$ws.Range("L9").Value = 0.861397757647421

$ws.Range("B10").Value = 0.672713993521846
$ws.Range("L10").Value = 0.527106989746053

$ws.Range("B11").Value = 0.531503319623507
$ws.Range("E11").Value = 0.566216773671888
$ws.Range("L11").Value = 0.547506346993746
